$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 226, shifting existing rows 226-280 down to 227-281
$ws.Rows.Item(226).Insert()

# Populate the newly inserted row 226 with the new data record
$ws.Cells.Item(226, 1).Value = 3
$ws.Cells.Item(226, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(226, 3).Value = "Coquimbo"
$ws.Cells.Item(226, 4).Value = 45135
$ws.Cells.Item(226, 5).Value = 5
$ws.Cells.Item(226, 6).Value = 100112026
$ws.Cells.Item(226, 7).Value = "Haba"
$ws.Cells.Item(226, 8).Value = "Sin especificar"
$ws.Cells.Item(226, 9).Value = "Primera"
$ws.Cells.Item(226, 10).Value = 55
$ws.Cells.Item(226, 11).Value = 15000
$ws.Cells.Item(226, 12).Value = 15000
$ws.Cells.Item(226, 13).Value = 15000
$ws.Cells.Item(226, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(226, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(226, 16).Value = 600
$ws.Cells.Item(226, 17).Value = 25
$ws.Cells.Item(226, 18).Value = "Hortaliza"
